$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 244. This shifts the former rows
# 244-350 down to 247-353, preserving their content and formatting.
$ws.Rows("244:246").Insert()

# Populate the 3 newly inserted rows with a new week's data
# (same market/region/category/variety/quality pattern as the old
# top-of-block rows, now shifted to 247-249), dated 44455, with the
# new Volumen / Precio minimo / Precio maximo / Precio promedio
# ponderado / Precio $/Kg values.

# Row 244: Larga vida / Primera
$ws.Range("A244").Value = 2
$ws.Range("B244").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C244").Value = "Coquimbo"
$ws.Range("D244").Value = 44455
$ws.Range("E244").Value = 4
$ws.Range("F244").Value = 100112020
$ws.Range("G244").Value = "Tomate"
$ws.Range("H244").Value = "Larga vida"
$ws.Range("I244").Value = "Primera"
$ws.Range("J244").Value = 2600
$ws.Range("K244").Value = 10500
$ws.Range("L244").Value = 11000
$ws.Range("M244").Value = 10750
$ws.Range("N244").Value = "$/bandeja 18 kilos"
$ws.Range("O244").Value = "Provincia de Limarí"
$ws.Range("P244").Value = 597
$ws.Range("Q244").Value = 18
$ws.Range("R244").Value = "Hortaliza"

# Row 245: Larga vida / Segunda
$ws.Range("A245").Value = 2
$ws.Range("B245").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C245").Value = "Coquimbo"
$ws.Range("D245").Value = 44455
$ws.Range("E245").Value = 4
$ws.Range("F245").Value = 100112020
$ws.Range("G245").Value = "Tomate"
$ws.Range("H245").Value = "Larga vida"
$ws.Range("I245").Value = "Segunda"
$ws.Range("J245").Value = 2200
$ws.Range("K245").Value = 8500
$ws.Range("L245").Value = 9000
$ws.Range("M245").Value = 8750
$ws.Range("N245").Value = "$/bandeja 18 kilos"
$ws.Range("O245").Value = "Provincia de Limarí"
$ws.Range("P245").Value = 486
$ws.Range("Q245").Value = 18
$ws.Range("R245").Value = "Hortaliza"

# Row 246: Larga vida / Tercera
$ws.Range("A246").Value = 2
$ws.Range("B246").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C246").Value = "Coquimbo"
$ws.Range("D246").Value = 44455
$ws.Range("E246").Value = 4
$ws.Range("F246").Value = 100112020
$ws.Range("G246").Value = "Tomate"
$ws.Range("H246").Value = "Larga vida"
$ws.Range("I246").Value = "Tercera"
$ws.Range("J246").Value = 1600
$ws.Range("K246").Value = 6500
$ws.Range("L246").Value = 7000
$ws.Range("M246").Value = 6750
$ws.Range("N246").Value = "$/bandeja 18 kilos"
$ws.Range("O246").Value = "Provincia de Limarí"
$ws.Range("P246").Value = 375
$ws.Range("Q246").Value = 18
$ws.Range("R246").Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of
# column D (style index 2 / "YYYY-MM-DD HH:MM:SS").
$ws.Range("D244:D246").NumberFormat = $ws.Range("D247").NumberFormat
